$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations - Market")

# Insert two new rows for the "nicotine" filter translation keys, right
# above the existing "market.filter.pgvg.filter.pgvg.label" row (row 46),
# keeping the sheet's alphabetical ordering by the Label column (B).
$ws.Rows("46:47").Insert()

$ws.Range("A46").Value2 = "cs"
$ws.Range("B46").Value2 = "market.filter.pgvg.filter.nicotine.label"
$ws.Range("C46").Value2 = "Nikotin"

$ws.Range("A47").Value2 = "cs"
$ws.Range("B47").Value2 = "market.filter.pgvg.filter.nicotine.reset"
$ws.Range("C47").Value2 = "Vše"

# Re-apply a sort (by Label, column B) over the whole table so the
# worksheet keeps its recorded sort state, matching what Excel stores
# after the user re-sorts the (now larger) table.
$rng = $ws.Range("A2:C77")
$key = $ws.Range("B2:B77")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($key)
$ws.Sort.SetRange($rng)
$ws.Sort.Header = -4142
$ws.Sort.Apply() | Out-Null

# Restore the active selection reported by the workbook after the edit.
$ws.Range("B63").Select() | Out-Null
